$wb = $excel.ActiveWorkbook

# Sheet "展览" (index 1 / sheet1.xml)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 1349
$ws1.Range("F6").Value = 7700
$ws1.Range("F14").Value = 5689
$ws1.Range("F16").Value = 2649
$ws1.Range("F24").Value = 3631
$ws1.Range("F29").Value = 3126
$ws1.Range("F30").Value = 58
$ws1.Range("F33").Value = 134
$ws1.Range("F34").Value = 330
$ws1.Range("F35").Value = 826
$ws1.Range("F39").Value = 2289
$ws1.Range("F43").Value = 3059

# Sheet "全部类型" (index 4 / sheet4.xml)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value = 1349
$ws4.Range("F6").Value = 7700
$ws4.Range("F13").Value = 5689
$ws4.Range("F15").Value = 2649
$ws4.Range("F25").Value = 3631
$ws4.Range("F30").Value = 3126
$ws4.Range("F31").Value = 58
$ws4.Range("F33").Value = 134
$ws4.Range("F34").Value = 330
$ws4.Range("F36").Value = 826
$ws4.Range("F41").Value = 2289
$ws4.Range("F45").Value = 3059
